$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.317.07"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.870.49"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.90"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4707"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2884"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06625"
$ws.Range("E9").Value = "  +1.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.74"
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08032"
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.47"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.869.75"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.149"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6870"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "271.22"
$ws.Range("E16").Value = "  -2.42%  "
$ws.Range("D17").Value = "30.315.92"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.15"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007740"
$ws.Range("E19").Value = "  +5.98%  "
$ws.Range("D21").Value = "2.116.85"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.309"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.225"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.422"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.77"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("E29").Value = "  -0.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09904"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.383"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.089"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  -0.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  +0.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7028"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.650"
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.313"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.85"
$ws.Range("E41").Value = "  -3.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.962"
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8443"
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4174"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.31"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.274"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.086"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "928.18"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.51"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05685"
$ws.Range("E51").Value = "  +0.63%  "
